$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected. Rather than Unprotect()/Protect() (which would
# rewrite <sheetProtection> with a brand-new hash algorithm and drop the
# original protection flags), temporarily unlock just the cells we need to
# touch, edit them, then relock them so the saved sheetProtection element
# stays exactly as it was.

$ws.Range("A10").Locked = $false
$ws.Range("D2:E7").Locked = $false

# Update the disclosure/date text: 2021-03-18 -> 2021-03-19
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-19 for illustrative purposes only and are subject to change."

# Update the Weight (D) / Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2442618712688534
$ws.Range("E2").Value = 0.001778726431874933

$ws.Range("D3").Value = 0.4971247477546838
$ws.Range("E3").Value = -0.005030743432084983

$ws.Range("D4").Value = 0.09786010203055576
$ws.Range("E4").Value = 0.007806080525883496

$ws.Range("D5").Value = 0.1016044850531913
$ws.Range("E5").Value = -0.002305475504323029

$ws.Range("D6").Value = 0.05914879389271555
$ws.Range("E6").Value = 0.002978235967926679

$ws.Range("D7").Value = 0.9999999999999998
$ws.Range("E7").Value = -0.001360615762235895

$ws.Range("A10").Locked = $true
$ws.Range("D2:E7").Locked = $true
